$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 82, shifting existing rows 82-186 down to 83-187
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with the new data
$ws.Cells.Item(82, 1).Value = 11
$ws.Cells.Item(82, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(82, 3).Value = "Bíobío"
$ws.Cells.Item(82, 4).Value = 44763
$ws.Cells.Item(82, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(82, 5).Value = 8
$ws.Cells.Item(82, 6).Value = "Fruta"
$ws.Cells.Item(82, 7).Value = 100108
$ws.Cells.Item(82, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(82, 9).Value = 100108005
$ws.Cells.Item(82, 10).Value = "Piña"
$ws.Cells.Item(82, 11).Value = "Caramelo"
$ws.Cells.Item(82, 12).Value = "Segunda"
$ws.Cells.Item(82, 13).Value = 200
$ws.Cells.Item(82, 14).Value = 18000
$ws.Cells.Item(82, 15).Value = 19000
$ws.Cells.Item(82, 16).Value = 18500
$ws.Cells.Item(82, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(82, 18).Value = "Ecuador"
$ws.Cells.Item(82, 19).Value = 1321
$ws.Cells.Item(82, 20).Value = 14
